# Applies the "Adicionando os últimos comandos" edit:
#  - marks several more commands as implemented in the Assembler ("F") column
#  - renames the "Bcc" instruction row label to "BCC"
#  - updates the view's top-left cell / active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# --- Instruction (C38): Bcc -> BCC -------------------------------------
$ws.Range("C38").Value = "BCC"

# --- Assembler column ("F") updates: no -> yes (fully or partially) ----
$ws.Range("F10").Value = "yes"
$ws.Range("F11").Value = "yes"
$ws.Range("F12").Value = "yes"
$ws.Range("F14").Value = "yes"
$ws.Range("F15").Value = "yes" + $nl + "???" + $nl + "no"
$ws.Range("F16").Value = "yes" + $nl + "??" + $nl + "no"
$ws.Range("F17").Value = "yes"
$ws.Range("F18").Value = "yes"
$ws.Range("F19").Value = "yes" + $nl + "no" + $nl + "??"
$ws.Range("F20").Value = "yes" + $nl + "no" + $nl + "??"
$ws.Range("F21").Value = "yes"
$ws.Range("F22").Value = "yes"
$ws.Range("F23").Value = "yes" + $nl + "yes" + $nl + "yes"
$ws.Range("F24").Value = "yes" + $nl + "yes" + $nl + "yes"
$ws.Range("F25").Value = "yes"
$ws.Range("F26").Value = "yes"
$ws.Range("F33").Value = "yes"
$ws.Range("F37").Value = "yes"
$ws.Range("F39").Value = "yes"
$ws.Range("F41").Value = "yes" + $nl + "yes"
$ws.Range("F42").Value = "yes" + $nl + "yes"
$ws.Range("F43").Value = "yes"

# --- Scroll / selection state -------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("D27").Select()
